$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.106720805168152
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.476559162139893
$ws.Range("D1").Value = 1.365588188171387
$ws.Range("E1").Value = 0.9881521463394165
